$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data (runs/balls/fours for rows 2 and 3) was swapped between
# the two Prabhsimran Singh rows. Force the target range to Text first so
# Excel doesn't reinterpret the numeric-looking strings as numbers (the
# original cells are stored as text), then restore default formatting so
# no extra style is left behind on the cells.
$rng = $ws.Range("C2:E3")
$rng.NumberFormat = "@"

$ws.Range("C2").Value = "11"
$ws.Range("D2").Value = "8"
$ws.Range("E2").Value = "2"

$ws.Range("C3").Value = "4"
$ws.Range("D3").Value = "7"
$ws.Range("E3").Value = "0"

$rng.ClearFormats()
